$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '28.019.50'
$ws.Range('E2').Value = '  -1.99%  '
$ws.Range('D3').Value = '1.830.76'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '324.63'
$ws.Range('E5').Value = '  -3.37%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4647'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '0.3866'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').Value = '0.07836'
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').Value = '0.9588'
$ws.Range('E10').Value = '  -2.40%  '
$ws.Range('D11').Value = '21.89'
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('D12').Value = '1.842.64'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').Value = '5.678'
$ws.Range('E13').Value = '  -2.90%  '
$ws.Range('D14').Value = '6.887'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '0.06857'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '88.24'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '0.000009916'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '28.024.58'
$ws.Range('E21').Value = '  -1.97%  '
$ws.Range('D22').Value = '5.292'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('E23').Value = '  -3.31%  '
$ws.Range('D24').Value = '2.088'
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').Value = '2.043.46'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '154.86'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = '19.12'
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').Value = '5.655'
$ws.Range('E28').Value = '  -6.35%  '
$ws.Range('E29').Value = '  -3.24%  '
$ws.Range('D30').Value = '118.38'
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('D31').Value = '0.09246'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D33').Value = '5.253'
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('E34').Value = '  -2.10%  '
$ws.Range('D35').Value = '3.308'
$ws.Range('E35').Value = '  -4.91%  '
$ws.Range('D36').Value = '0.05840'
$ws.Range('E36').Value = '  -5.00%  '
$ws.Range('D37').Value = '0.02121'
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('D38').Value = '1.142'
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').Value = '7.717'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('D40').Value = '0.5586'
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('D42').Value = '0.1755'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('D43').Value = '0.07258'
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('D44').Value = '11.59'
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('D45').Value = '0.5258'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('D46').Value = '1.145'
$ws.Range('E46').Value = '  -8.37%  '
$ws.Range('D47').Value = '2.094'
$ws.Range('E47').Value = '  -11.35%  '
$ws.Range('D48').Value = '1.819'
$ws.Range('E48').Value = '  -4.76%  '
$ws.Range('D49').Value = '112.83'
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').Value = '2.322'
$ws.Range('E51').Value = '  +0.44%  '
